$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 20003020
$ws.Range("I32").Value = 1750
$ws.Range("K32").Value = 1750
$ws.Range("M32").Value = -1424

$ws.Range("H40").Value = 7128.5713
$ws.Range("I40").Value = 9900
$ws.Range("J40").Value = 6666.6665
$ws.Range("K40").Value = 9900
$ws.Range("L40").Value = 6666.6665
$ws.Range("M40").Value = -9725
$ws.Range("N40").Value = -7016.6665

$ws.Range("H104").Value = 147.66667
$ws.Range("I104").Value = 147.66667
$ws.Range("K104").Value = 443.00001
$ws.Range("M104").Value = 1303.99999

$ws.Range("H109").Value = 39500
$ws.Range("J109").Value = 39500
$ws.Range("L109").Value = 39500
$ws.Range("N109").Value = -42274

$ws.Range("H127").Value = 5121.846
$ws.Range("I127").Value = 2069.3333
$ws.Range("K127").Value = 6207.999899999999
$ws.Range("M127").Value = -1247.999899999999

$ws.Range("H129").Value = 168554.58
$ws.Range("I129").Value = 183695.9
$ws.Range("K129").Value = 551087.7
$ws.Range("M129").Value = -546087.7

$ws.Range("H132").Value = 14190.878
$ws.Range("I132").Value = 1652.3143
$ws.Range("J132").Value = 87332.5
$ws.Range("K132").Value = 4956.9429
$ws.Range("L132").Value = 261997.5
$ws.Range("M132").Value = -2426.9429
$ws.Range("N132").Value = -267057.5

$ws.Range("H138").Value = 2100.3215
$ws.Range("I138").Value = 1498.1364
$ws.Range("J138").Value = 4308.3335
$ws.Range("K138").Value = 4494.4092
$ws.Range("L138").Value = 12925.0005
$ws.Range("M138").Value = 645.5907999999999
$ws.Range("N138").Value = -23205.0005

$ws.Range("H141").Value = 8205.866
$ws.Range("I141").Value = 7655.909
$ws.Range("K141").Value = 22967.727
$ws.Range("M141").Value = -17787.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H48").Value = 380000
$ws.Range("J48").Value = 380000
$ws.Range("L48").Value = 380000
$ws.Range("N48").Value = -380768

$ws.Range("H74").Value = 2617.72
$ws.Range("I74").Value = 2226.7917
$ws.Range("K74").Value = 2226.7917
$ws.Range("M74").Value = -1352.7917

$ws.Range("H77").Value = 2617.72
$ws.Range("I77").Value = 2226.7917
$ws.Range("K77").Value = 11133.9585
$ws.Range("M77").Value = -6765.958500000001

$ws.Range("H118").Value = 39833.332
$ws.Range("J118").Value = 39833.332
$ws.Range("L118").Value = 39833.332
$ws.Range("N118").Value = -43147.332

$ws.Range("H122").Value = 4739.3447
$ws.Range("I122").Value = 4202.4707
$ws.Range("J122").Value = 5499.9165
$ws.Range("K122").Value = 12607.4121
$ws.Range("L122").Value = 16499.7495
$ws.Range("M122").Value = -10157.4121
$ws.Range("N122").Value = -21399.7495

$ws.Range("H132").Value = 1407.9445
$ws.Range("I132").Value = 1334.25
$ws.Range("J132").Value = 1997.5
$ws.Range("K132").Value = 4002.75
$ws.Range("L132").Value = 5992.5
$ws.Range("M132").Value = -1472.75
$ws.Range("N132").Value = -11052.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2478.8
$ws.Range("I105").Value = 1490.4286
$ws.Range("K105").Value = 1490.4286
$ws.Range("M105").Value = 256.5714

$ws.Range("H107").Value = 11661.4
$ws.Range("I107").Value = 10499.637
$ws.Range("K107").Value = 10499.637
$ws.Range("M107").Value = -8579.637000000001

$ws.Range("H134").Value = 2448.0667
$ws.Range("I134").Value = 2029.8182
$ws.Range("K134").Value = 6089.4546
$ws.Range("M134").Value = -3554.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3361.182
$ws.Range("I31").Value = 2320
$ws.Range("J31").Value = 4228.8335
$ws.Range("K31").Value = 2320
$ws.Range("L31").Value = 4228.8335
$ws.Range("M31").Value = -2025
$ws.Range("N31").Value = -4818.8335

$ws.Range("H34").Value = 3361.182
$ws.Range("I34").Value = 2320
$ws.Range("J34").Value = 4228.8335
$ws.Range("K34").Value = 2320
$ws.Range("L34").Value = 4228.8335
$ws.Range("M34").Value = -2118
$ws.Range("N34").Value = -4632.8335

$ws.Range("H99").Value = 2812.375

$ws.Range("H107").Value = 9009.25
$ws.Range("J107").Value = 17770
$ws.Range("L107").Value = 17770
$ws.Range("N107").Value = -21610

$ws.Range("H122").Value = 2808.4783
$ws.Range("I122").Value = 2142.5
$ws.Range("K122").Value = 6427.5
$ws.Range("M122").Value = -3977.5

$ws.Range("H126").Value = 2812.375

$ws.Range("H132").Value = 2611.5833
$ws.Range("I132").Value = 2289
$ws.Range("K132").Value = 6867
$ws.Range("M132").Value = -4337

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 49219428
$ws.Range("I4").Value = 51680148
$ws.Range("K4").Value = 155040444
$ws.Range("M4").Value = -155040332

$ws.Range("H121").Value = 1521.7778
$ws.Range("J121").Value = 1671.5714
$ws.Range("L121").Value = 5014.7142
$ws.Range("N121").Value = -7634.7142

$ws.Range("H129").Value = 1727.7858
$ws.Range("J129").Value = 2380.4
$ws.Range("L129").Value = 7141.200000000001
$ws.Range("N129").Value = -17141.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 19000
$ws.Range("I43").Value = 19000
$ws.Range("K43").Value = 19000
$ws.Range("M43").Value = -18849

$ws.Range("H80").Value = 89619.84
$ws.Range("I80").Value = 114006
$ws.Range("J80").Value = 8332.666999999999
$ws.Range("K80").Value = 114006
$ws.Range("L80").Value = 8332.666999999999
$ws.Range("M80").Value = -113008
$ws.Range("N80").Value = -10328.667

$ws.Range("H83").Value = 89619.84
$ws.Range("I83").Value = 114006
$ws.Range("J83").Value = 8332.666999999999
$ws.Range("K83").Value = 570030
$ws.Range("L83").Value = 41663.335
$ws.Range("M83").Value = -565038
$ws.Range("N83").Value = -51647.335

$ws.Range("H97").Value = 456.42856
$ws.Range("I97").Value = 394.14285
$ws.Range("J97").Value = 518.7143
$ws.Range("K97").Value = 394.14285
$ws.Range("L97").Value = 518.7143
$ws.Range("M97").Value = 101.85715
$ws.Range("N97").Value = -1510.7143

$ws.Range("H102").Value = 5331.4287
$ws.Range("I102").Value = 4787.8945
$ws.Range("K102").Value = 4787.8945
$ws.Range("M102").Value = -3165.8945

$ws.Range("H122").Value = 5705.3335
$ws.Range("I122").Value = 5015.077
$ws.Range("K122").Value = 15045.231
$ws.Range("M122").Value = -12595.231

$ws.Range("H132").Value = 6451.553
$ws.Range("I132").Value = 5834.5674
$ws.Range("K132").Value = 17503.7022
$ws.Range("M132").Value = -14973.7022

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3129.7144
$ws.Range("I7").Value = 2651.3333
$ws.Range("K7").Value = 2651.3333
$ws.Range("M7").Value = -2539.3333

$ws.Range("H25").Value = 230002.33
$ws.Range("I25").Value = 230002.33
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 230002.33
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -229772.33
$ws.Range("N25").ClearContents()

$ws.Range("H55").Value = 387.67856
$ws.Range("I55").Value = 369.94116
$ws.Range("J55").Value = 415.0909
$ws.Range("K55").Value = 369.94116
$ws.Range("L55").Value = 415.0909
$ws.Range("M55").Value = -196.94116
$ws.Range("N55").Value = -761.0908999999999

$ws.Range("H68").Value = 3750.077
$ws.Range("I68").Value = 3616
$ws.Range("K68").Value = 3616
$ws.Range("M68").Value = -2867

$ws.Range("H71").Value = 3750.077
$ws.Range("I71").Value = 3616
$ws.Range("K71").Value = 18080
$ws.Range("M71").Value = -14336

$ws.Range("H82").Value = 31251556
$ws.Range("I82").Value = 45456140
$ws.Range("J82").Value = 1471.6
$ws.Range("K82").Value = 45456140
$ws.Range("L82").Value = 1471.6
$ws.Range("M82").Value = -45455779
$ws.Range("N82").Value = -2193.6

$ws.Range("H85").Value = 31251556
$ws.Range("I85").Value = 45456140
$ws.Range("J85").Value = 1471.6
$ws.Range("K85").Value = 45456140
$ws.Range("L85").Value = 1471.6
$ws.Range("M85").Value = -45454892
$ws.Range("N85").Value = -3967.6

$ws.Range("H122").Value = 7024.75
$ws.Range("I122").Value = 6999.6665
$ws.Range("K122").Value = 20998.9995
$ws.Range("M122").Value = -18548.9995

$ws.Range("H126").Value = 3129.7144
$ws.Range("I126").Value = 2651.3333
$ws.Range("K126").Value = 7953.999899999999
$ws.Range("M126").Value = -5483.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 15000
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H109").Value = 17944.445
$ws.Range("J109").Value = 17944.445
$ws.Range("L109").Value = 17944.445
$ws.Range("N109").Value = -20718.445

$ws.Range("H115").Value = 29833.334
$ws.Range("J115").Value = 29833.334
$ws.Range("L115").Value = 29833.334
$ws.Range("N115").Value = -32967.334

$ws.Range("H122").Value = 1925.4667
$ws.Range("I122").Value = 1480.1818
$ws.Range("K122").Value = 4440.5454
$ws.Range("M122").Value = -1990.5454

$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530

$ws.Range("H136").Value = 55002.223
$ws.Range("I136").Value = 61018.125
$ws.Range("K136").Value = 183054.375
$ws.Range("M136").Value = -180504.375
